# Add a new "expectedResultForReg" column (D) to the data sheet, used by the
# duplicate-registration DDT case: every existing data row gets the literal
# "duplicate" expected result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 31
$headerCol = "D1"
$dataRange = "D2:D" + $lastRow

# Header cell D1 - reuse the same look as the other header cells (A1:C1).
$ws.Range("A1").Copy()
$ws.Range($headerCol).PasteSpecial(-4122)
$ws.Range($headerCol).Value = "expectedResultForReg"

# Data cells D2:D31 - reuse the same look as the other data cells (A2:C31).
$ws.Range("A2").Copy()
$ws.Range($dataRange).PasteSpecial(-4122)
$ws.Range($dataRange).Value = "duplicate"

# Size the new column and move the active selection onto it, like Excel does
# right after typing a new header.
$ws.Columns.Item(4).ColumnWidth = 20.86
$ws.Range("D1").Select() | Out-Null

# Keep the rest of the sheet's print setup explicit (portrait).
$ws.PageSetup.Orientation = 1
